$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of rows 5, 6 and 7 (new row5 = old row7,
# new row6 = old row5, new row7 = old row6), rounds the Ost/Nord
# (Q/R) coordinates to whole numbers, drops the Starttid/Sluttid
# (Z/AB) values on all three rows, and moves the "Publik kommentar"
# (AC) text along with its row.

# --- Row 5 (becomes what used to be row 7) ---
$ws.Range("A5").Value = 111817611
$ws.Range("B5").Value = 89416
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 1205
$ws.Range("F5").Value = "Stor aspticka"
$ws.Range("G5").Value = "Phellinus populicola"
$ws.Range("H5").Value = "Niemelä"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "1"
$ws.Range("Q5").Value = 578480
$ws.Range("R5").Value = 6398700
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("AC5").Value = "Relativt tunn asp."

# --- Row 6 (becomes what used to be row 5) ---
$ws.Range("A6").Value = 111817654
$ws.Range("B6").Value = 89363
$ws.Range("E6").Value = 5445
$ws.Range("F6").Value = "Ekticka"
$ws.Range("G6").Value = "Fomitiporia robusta"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "5"
$ws.Range("Q6").Value = 578450
$ws.Range("R6").Value = 6398641
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AC6").Value = "På relativt tunn ek."

# --- Row 7 (becomes what used to be row 6) ---
$ws.Range("A7").Value = 111817582
$ws.Range("B7").Value = 88283
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 655
$ws.Range("F7").Value = "Oxtungssvamp"
$ws.Range("G7").Value = "Fistulina hepatica"
$ws.Range("H7").Value = "(Schaeff.) With., nom sanct."
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "2"
$ws.Range("Q7").Value = 578499
$ws.Range("R7").Value = 6398731
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").ClearContents()
